$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FAIL")
$ws.Activate()
